# BancoDados.xlsx - "Implement product deletion functionality" data update
#
# - Row 2 (previously Caio Cesar's entry) becomes Caio Gonzaga's new entry
# - Row 3 (Bruno Auzier's entry) gets new defect/description values
# - Row 4 (Maycon Nascimento's entry) gets new defect/description values
# - Rows 5 (Victor Barbosa) and 6 (old placeholder/test row) are removed
#   (this is the "DeleteProduct" functionality mentioned in the commit)
# - Selection moves to A7, the first empty row below the remaining table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Turno" column (B) holds small integers stored as TEXT (shared string),
# not numbers. Typing a bare numeric string via .Value would make Excel
# coerce it to a real number, so instead we copy the already-text-typed
# value from another cell in the same column (values-only paste) - this
# keeps the cell as text without touching styles.xml at all.
# Do this BEFORE overwriting the source cells below.
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4163)   # xlPasteValues -> B4 becomes "2" (text)

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4163)   # xlPasteValues -> B3 becomes "1" (text)
# B2 is already the text "1", so it needs no change.
$excel.CutCopyMode = 0

# --- Row 2: Caio Gonzaga ---
$ws.Range("A2").Value = "Caio Gonzaga"
$ws.Range("C2").Value = "FE01"
$ws.Range("D2").Value = "FE FRESH"
$ws.Range("E2").Value = "ZF5254T7ZL"
$ws.Range("F2").Value = "LAMU"
$ws.Range("G2").Value = "Displaced weld deposition"
$ws.Range("H2").Value = "Solder"
$ws.Range("I2").Value = "Shifted"
$ws.Range("J2").Value = "Teste.`nTeste."
$ws.Range("K2").Value = "3/17/2025 3:50 PM"

# --- Row 3: Bruno Auzier ---
$ws.Range("A3").Value = "Bruno Auzier"
$ws.Range("C3").Value = "FE05"
$ws.Range("D3").Value = "FE FRESH"
$ws.Range("E3").Value = "ZF5254T7ZL"
$ws.Range("F3").Value = "LAMU"
$ws.Range("G3").Value = "Lack of PM, obstructed, damaged, worn, life time, etc."
$ws.Range("H3").Value = "Maintenance"
$ws.Range("I3").Value = "Nozzle"
$ws.Range("J3").Value = "Teste.`nTese."
$ws.Range("K3").Value = "3/17/2025 3:51 PM"

# --- Row 4: Maycon Nascimento ---
$ws.Range("A4").Value = "Maycon Nascimento"
$ws.Range("C4").Value = "BE07"
$ws.Range("D4").Value = "BE FRESH"
$ws.Range("E4").Value = "ZF5254T7XJ"
$ws.Range("F4").Value = "MANILA"
$ws.Range("G4").Value = "Engineering activities (Process engineering, maintenance, testing, product, etc.)"
$ws.Range("H4").Value = "Engineering"
$ws.Range("I4").Value = "Beta process"
$ws.Range("J4").Value = "Teste.`nTeste."
$ws.Range("K4").Value = "3/17/2025 3:53 PM"

# The multi-line descriptions above trigger Excel's automatic row-height
# fit; re-run AutoFit so the rows end up back at the default height (no
# explicit row height override), matching the original workbook.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()

# --- Remove the two trailing rows (Victor Barbosa + test/placeholder row) ---
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# Move the active selection to the first free row below the data
$ws.Range("A7").Select()
